$d = $word.ActiveDocument

# Replace the text of the first paragraph.
$d.Content.Find.Execute(
    "Video provides a powerful way to help you prove your point. When you click Online Video, you can paste in the embed code for the video you want to add. You can also type a keyword to search online for the video that best fits your document.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This is after the modification.", 2)

# Delete paragraphs 2 through 5 (the four "feature" paragraphs), leaving the
# trailing empty paragraph intact.
$p2 = $d.Paragraphs.Item(2)
$p5 = $d.Paragraphs.Item(5)
$r = $d.Range($p2.Range.Start, $p5.Range.End)
$r.Delete()
